$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 (legmaxROM header values)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 (meanEMG) - B2 cleared entirely, C2:E2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 0.90008561529877906
$ws.Range("D2").Value = 0.14043214109872218
$ws.Range("E2").Value = 3.2613619788749233

# Update row 3 (meanEMG) - B3:E3 updated
$ws.Range("B3").Value = 0.22171485292124826
$ws.Range("C3").Value = 1.1008864504096048
$ws.Range("D3").Value = 0.26917189315422441
$ws.Range("E3").Value = 0.93931927950181182

# Update selection to match new reduced range of interest
$ws.Range("B1:E3").Select()
